$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 101, shifting existing rows 101-120 down to 102-121.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new weekly data record.
$ws.Range("A101").Value = 6
$ws.Range("B101").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C101").Value = "Metropolitana"
$ws.Range("D101").Value = 45015
$ws.Range("E101").Value = 13
$ws.Range("F101").Value = 100114007
$ws.Range("G101").Value = "Jengibre"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 250
$ws.Range("K101").Value = 15000
$ws.Range("L101").Value = 16000
$ws.Range("M101").Value = 15400
$ws.Range("N101").Value = "$/caja 13 kilos"
$ws.Range("O101").Value = "Perú"
$ws.Range("P101").Value = 1185
$ws.Range("Q101").Value = 13
$ws.Range("R101").Value = "Hortaliza"
